$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-case text cells affected by the unit test changes
$ws.Range("E17").Value = "from bank_account.bank_account import BankAccount                                    import unittest "
$ws.Range("E18").Value = "from bank_account.bank_account import BankAccount                                    import unittest "
$ws.Range("E19").Value = "from bank_account.bank_account import BankAccount                                    import unittest"
$ws.Range("E20").Value = "from bank_account.bank_account import BankAccount                                    import unittest"
$ws.Range("E21").Value = "from bank_account.bank_account import BankAccount                                    import unittest "
$ws.Range("F17").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                                       actual = BankAccount.deosit(self.bank_account, 21)                      self.assertIsNone(actual)"
$ws.Range("F16").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                              actual=   BankAccount.update_balance(self.bank_account, `"Ridham`")     self.assertIsNone(actual)                               "
$ws.Range("F14").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                                                  actual=BankAccount.update_balance(self.bank_account, 21)                                    self.assertIsNone"
$ws.Range("F15").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                                                  actual=BankAccount.update_balance(self.bank_account, -21)                              self.assertIsNone"
$ws.Range("F19").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                        actual = BankAccount.deosit(self.bank_account, 21)                      self.assertIsNone(actual)"
$ws.Range("F18").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                        actual = BankAccount.deosit(self.bank_account, -21)                      self.assertIsNone(actual)"
$ws.Range("F20").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                        actual = BankAccount.deosit(self.bank_account,- 21)                      self.assertIsNone(actual)"
$ws.Range("F21").Value = "bank_account = BankAccount(20021, 123, 1123.123)                                                        actual = BankAccount.deosit(self.bank_account, 3221)                      self.assertIsNone(actual)"

# Row heights grew to fit the longer test code snippets
$ws.Rows("14:14").RowHeight = 88.2
$ws.Rows("15:15").RowHeight = 93.6
$ws.Rows("16:16").RowHeight = 105
$ws.Rows("17:17").RowHeight = 105.6
$ws.Rows("18:18").RowHeight = 107.4
$ws.Rows("19:19").RowHeight = 99.6
$ws.Rows("20:20").RowHeight = 109.2
$ws.Rows("21:21").RowHeight = 102.6

# Restore the selection/scroll position as last edited
[void]$ws.Range("F21").Select()
